$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H105").Value = 78671
$ws.Range("J105").Value = 78671
$ws.Range("L105").Value = 78671
$ws.Range("N105").Value = -85659
$ws.Range("H113").Value = 2306
$ws.Range("I113").Value = 2257.8572
$ws.Range("J113").Value = 2348.125
$ws.Range("K113").Value = 2257.8572
$ws.Range("L113").Value = 2348.125
$ws.Range("M113").Value = 996.1428000000001
$ws.Range("N113").Value = -8856.125
$ws.Range("H132").Value = 2262.75
$ws.Range("I132").Value = 2321.9524
$ws.Range("J132").Value = 2014.1
$ws.Range("K132").Value = 6965.8572
$ws.Range("L132").Value = 6042.299999999999
$ws.Range("M132").Value = -4435.8572
$ws.Range("N132").Value = -11102.3
$ws.Range("H137").Value = 1853084.2
$ws.Range("I137").Value = 3205923.5
$ws.Range("J137").Value = 1830.7368
$ws.Range("K137").Value = 9617770.5
$ws.Range("L137").Value = 5492.2104
$ws.Range("M137").Value = -9615220.5
$ws.Range("N137").Value = -10592.2104
$ws.Range("H141").Value = 2800.1292
$ws.Range("I141").Value = 1246.625
$ws.Range("K141").Value = 3739.875
$ws.Range("M141").Value = 1440.125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6864152
$ws.Range("I32").Value = 8077573
$ws.Range("J32").Value = 24868.637
$ws.Range("K32").Value = 8077573
$ws.Range("L32").Value = 24868.637
$ws.Range("M32").Value = -8077286
$ws.Range("N32").Value = -25442.637
$ws.Range("H61").Value = 15876317
$ws.Range("I61").Value = 20835416
$ws.Range("K61").Value = 20835416
$ws.Range("M61").Value = -20835204
$ws.Range("H74").Value = 13516072
$ws.Range("I74").Value = 2080.7896
$ws.Range("J74").Value = 27780840
$ws.Range("K74").Value = 2080.7896
$ws.Range("L74").Value = 27780840
$ws.Range("M74").Value = -1206.7896
$ws.Range("N74").Value = -27782588
$ws.Range("H77").Value = 13516072
$ws.Range("I77").Value = 2080.7896
$ws.Range("J77").Value = 27780840
$ws.Range("K77").Value = 10403.948
$ws.Range("L77").Value = 138904200
$ws.Range("M77").Value = -6035.948
$ws.Range("N77").Value = -138912936
$ws.Range("H132").Value = 1263896.8
$ws.Range("I132").Value = 2264.634
$ws.Range("K132").Value = 6793.902
$ws.Range("M132").Value = -4263.902
$ws.Range("H136").Value = 15876317
$ws.Range("I136").Value = 20835416
$ws.Range("K136").Value = 62506248
$ws.Range("M136").Value = -62503698
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5844.377
$ws.Range("I31").Value = 1880.7894
$ws.Range("J31").Value = 7350.54
$ws.Range("K31").Value = 1880.7894
$ws.Range("L31").Value = 7350.54
$ws.Range("M31").Value = -1585.7894
$ws.Range("N31").Value = -7940.54
$ws.Range("H34").Value = 5844.377
$ws.Range("I34").Value = 1880.7894
$ws.Range("J34").Value = 7350.54
$ws.Range("K34").Value = 1880.7894
$ws.Range("L34").Value = 7350.54
$ws.Range("M34").Value = -1678.7894
$ws.Range("N34").Value = -7754.54
$ws.Range("H58").Value = 1101.7273
$ws.Range("I58").Value = 813.06665
$ws.Range("J58").Value = 1720.2858
$ws.Range("K58").Value = 813.06665
$ws.Range("L58").Value = 1720.2858
$ws.Range("M58").Value = -610.06665
$ws.Range("N58").Value = -2126.2858
$ws.Range("H62").Value = 3782.36
$ws.Range("I62").Value = 4024.75
$ws.Range("J62").Value = 3558.6155
$ws.Range("K62").Value = 4024.75
$ws.Range("L62").Value = 3558.6155
$ws.Range("M62").Value = -3400.75
$ws.Range("N62").Value = -4806.6155
$ws.Range("H65").Value = 3782.36
$ws.Range("I65").Value = 4024.75
$ws.Range("J65").Value = 3558.6155
$ws.Range("K65").Value = 20123.75
$ws.Range("L65").Value = 17793.0775
$ws.Range("M65").Value = -17003.75
$ws.Range("N65").Value = -24033.0775
$ws.Range("H122").Value = 2284.5186
$ws.Range("I122").Value = 2112.1333
$ws.Range("K122").Value = 6336.3999
$ws.Range("M122").Value = -3886.3999
$ws.Range("H132").Value = 22224626
$ws.Range("I132").Value = 26317904
$ws.Range("J132").Value = 15154417
$ws.Range("K132").Value = 78953712
$ws.Range("L132").Value = 45463251
$ws.Range("M132").Value = -78951182
$ws.Range("N132").Value = -45468311
$ws.Range("H136").Value = 1101.7273
$ws.Range("I136").Value = 813.06665
$ws.Range("J136").Value = 1720.2858
$ws.Range("K136").Value = 2439.19995
$ws.Range("L136").Value = 5160.857400000001
$ws.Range("M136").Value = 110.8000499999998
$ws.Range("N136").Value = -10260.8574
$ws.Range("H138").Value = 40640
$ws.Range("J138").Value = 40640
$ws.Range("L138").Value = 40640
$ws.Range("N138").Value = -50920
$ws.Range("H140").Value = 59449.9
$ws.Range("J140").Value = 59449.9
$ws.Range("L140").Value = 59449.9
$ws.Range("N140").Value = -69809.89999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 568.03845
$ws.Range("I113").Value = 776
$ws.Range("J113").Value = 518.5238000000001
$ws.Range("K113").Value = 2328
$ws.Range("L113").Value = 1555.5714
$ws.Range("M113").Value = -158
$ws.Range("N113").Value = -5895.571400000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1880821
$ws.Range("I80").Value = 2251026.2
$ws.Range("J80").Value = 400000
$ws.Range("K80").Value = 2251026.2
$ws.Range("L80").Value = 400000
$ws.Range("M80").Value = -2250028.2
$ws.Range("N80").Value = -401996
$ws.Range("H83").Value = 1880821
$ws.Range("I83").Value = 2251026.2
$ws.Range("J83").Value = 400000
$ws.Range("K83").Value = 11255131
$ws.Range("L83").Value = 2000000
$ws.Range("M83").Value = -11250139
$ws.Range("N83").Value = -2009984
$ws.Range("H132").Value = 31255932
$ws.Range("I132").Value = 52639064
$ws.Range("J132").Value = 3664.3845
$ws.Range("K132").Value = 157917192
$ws.Range("L132").Value = 10993.1535
$ws.Range("M132").Value = -157914662
$ws.Range("N132").Value = -16053.1535
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 19250
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 19250
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 19250
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -20270
$ws.Range("H58").Value = 30000
$ws.Range("J58").Value = 30000
$ws.Range("L58").Value = 30000
$ws.Range("N58").Value = -30616
$ws.Range("H132").Value = 6785777.5
$ws.Range("I132").Value = 3078.8215
$ws.Range("J132").Value = 19446814
$ws.Range("K132").Value = 9236.4645
$ws.Range("L132").Value = 58340442
$ws.Range("M132").Value = -6706.4645
$ws.Range("N132").Value = -58345502
